# Update pricing/profit figures across the Leve profit sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR) to reflect refreshed market-board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 64.75
$ws.Range("I2").Value = 64.75
$ws.Range("K2").Value = 64.75
$ws.Range("M2").Value = 48.25
$ws.Range("H40").Value = 2000
$ws.Range("J40").Value = 2000
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2350
$ws.Range("H123").Value = 67593.336
$ws.Range("J123").Value = 67593.336
$ws.Range("L123").Value = 67593.336
$ws.Range("N123").Value = -77393.336
$ws.Range("H127").Value = 1000000
$ws.Range("I127").Value = 1000000
$ws.Range("K127").Value = 3000000
$ws.Range("M127").Value = -2995040
$ws.Range("H132").Value = 12282.954
$ws.Range("I132").Value = 13286.25
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 39858.75
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -37328.75
$ws.Range("N132").Value = -11810
$ws.Range("H137").Value = 12829748
$ws.Range("I137").Value = 25003108
$ws.Range("J137").Value = 15686.053
$ws.Range("K137").Value = 75009324
$ws.Range("L137").Value = 47058.159
$ws.Range("M137").Value = -75006774
$ws.Range("N137").Value = -52158.159
$ws.Range("H138").Value = 7709.886
$ws.Range("I138").Value = 6913.5
$ws.Range("K138").Value = 20740.5
$ws.Range("M138").Value = -15600.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 830203.6
$ws.Range("I32").Value = 925272
$ws.Range("K32").Value = 925272
$ws.Range("M32").Value = -924985
$ws.Range("H74").Value = 515329.38
$ws.Range("I74").Value = 574434.4
$ws.Range("J74").Value = 16220.444
$ws.Range("K74").Value = 574434.4
$ws.Range("L74").Value = 16220.444
$ws.Range("M74").Value = -573560.4
$ws.Range("N74").Value = -17968.444
$ws.Range("H77").Value = 515329.38
$ws.Range("I77").Value = 574434.4
$ws.Range("J77").Value = 16220.444
$ws.Range("K77").Value = 2872172
$ws.Range("L77").Value = 81102.22
$ws.Range("M77").Value = -2867804
$ws.Range("N77").Value = -89838.22
$ws.Range("H122").Value = 1540.5294
$ws.Range("I122").Value = 1486.1538
$ws.Range("J122").Value = 1717.25
$ws.Range("K122").Value = 4458.4614
$ws.Range("L122").Value = 5151.75
$ws.Range("M122").Value = -2008.4614
$ws.Range("N122").Value = -10051.75
$ws.Range("H132").Value = 4832.8433
$ws.Range("I132").Value = 3659.9583
$ws.Range("J132").Value = 5875.407
$ws.Range("K132").Value = 10979.8749
$ws.Range("L132").Value = 17626.221
$ws.Range("M132").Value = -8449.874899999999
$ws.Range("N132").Value = -22686.221
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 9045.857
$ws.Range("I99").Value = 9918
$ws.Range("J99").Value = 6865.5
$ws.Range("K99").Value = 9918
$ws.Range("L99").Value = 6865.5
$ws.Range("M99").Value = -8420
$ws.Range("N99").Value = -9861.5
$ws.Range("H132").Value = 383876.66
$ws.Range("J132").Value = 383876.66
$ws.Range("L132").Value = 383876.66
$ws.Range("N132").Value = -393996.66
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 716072.6
$ws.Range("I31").Value = 905240.3
$ws.Range("J31").Value = 6693.75
$ws.Range("K31").Value = 905240.3
$ws.Range("L31").Value = 6693.75
$ws.Range("M31").Value = -904945.3
$ws.Range("N31").Value = -7283.75
$ws.Range("H34").Value = 716072.6
$ws.Range("I34").Value = 905240.3
$ws.Range("J34").Value = 6693.75
$ws.Range("K34").Value = 905240.3
$ws.Range("L34").Value = 6693.75
$ws.Range("M34").Value = -905038.3
$ws.Range("N34").Value = -7097.75
$ws.Range("H105").Value = 20220.188
$ws.Range("I105").Value = 20220.188
$ws.Range("K105").Value = 20220.188
$ws.Range("M105").Value = -18473.188
$ws.Range("H134").Value = 4483.271
$ws.Range("I134").Value = 2117.1667
$ws.Range("J134").Value = 5902.933
$ws.Range("K134").Value = 6351.500100000001
$ws.Range("L134").Value = 17708.799
$ws.Range("M134").Value = -3816.500100000001
$ws.Range("N134").Value = -22778.799
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5145.295
$ws.Range("J68").Value = 5196.8306
$ws.Range("L68").Value = 15590.4918
$ws.Range("N68").Value = -17212.4918
$ws.Range("H71").Value = 5145.295
$ws.Range("J71").Value = 5196.8306
$ws.Range("L71").Value = 46771.4754
$ws.Range("N71").Value = -54883.4754
$ws.Range("H107").Value = 6508.154
$ws.Range("I107").Value = 4320.6
$ws.Range("K107").Value = 12961.8
$ws.Range("M107").Value = -11041.8
$ws.Range("H122").Value = 621444.0600000001
$ws.Range("I122").Value = 1467029.1
$ws.Range("J122").Value = 1348.3334
$ws.Range("K122").Value = 13203261.9
$ws.Range("L122").Value = 12135.0006
$ws.Range("M122").Value = -13200811.9
$ws.Range("N122").Value = -17035.0006
$ws.Range("H131").Value = 3704.3333
$ws.Range("J131").Value = 4204.4614
$ws.Range("L131").Value = 12613.3842
$ws.Range("N131").Value = -22693.3842
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 34999
$ws.Range("J68").Value = 34999
$ws.Range("L68").Value = 34999
$ws.Range("N68").Value = -36621
$ws.Range("H71").Value = 34999
$ws.Range("J71").Value = 34999
$ws.Range("L71").Value = 104997
$ws.Range("N71").Value = -113109
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H93").Value = 33000
$ws.Range("J93").Value = 33000
$ws.Range("L93").Value = 33000
$ws.Range("M93").Value = -36744
$ws.Range("H122").Value = 101700.3
$ws.Range("I122").Value = 101700.3
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 305100.9
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -302650.9
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 10427.116
$ws.Range("I132").Value = 17667.715
$ws.Range("J132").Value = 6931.6553
$ws.Range("K132").Value = 53003.145
$ws.Range("L132").Value = 20794.9659
$ws.Range("M132").Value = -50473.145
$ws.Range("N132").Value = -25854.9659
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 300.275
$ws.Range("I55").Value = 113.86957
$ws.Range("J55").Value = 552.4706
$ws.Range("K55").Value = 113.86957
$ws.Range("L55").Value = 552.4706
$ws.Range("M55").Value = 59.13043
$ws.Range("N55").Value = -898.4706
$ws.Range("H132").Value = 2781282.8
$ws.Range("I132").Value = 2979670.5
$ws.Range("J132").Value = 3852.5
$ws.Range("K132").Value = 8939011.5
$ws.Range("L132").Value = 11557.5
$ws.Range("M132").Value = -8936481.5
$ws.Range("N132").Value = -16617.5
$ws.Range("H136").Value = 14708422
$ws.Range("I136").Value = 8930948
$ws.Range("J136").Value = 41669970
$ws.Range("K136").Value = 26792844
$ws.Range("L136").Value = 125009910
$ws.Range("M136").Value = -26790294
$ws.Range("N136").Value = -125015010
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 44153
$ws.Range("I45").Value = 44994
$ws.Range("J45").Value = 43984.8
$ws.Range("K45").Value = 44994
$ws.Range("L45").Value = 43984.8
$ws.Range("M45").Value = -44503
$ws.Range("N45").Value = -44966.8
$ws.Range("H100").Value = 685
$ws.Range("I100").Value = 711.4286
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 1422.8572
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -881.8571999999999
$ws.Range("N100").Value = -2082
$ws.Range("H107").Value = 3849.1
$ws.Range("I107").Value = 1416.5
$ws.Range("J107").Value = 7498
$ws.Range("K107").Value = 4249.5
$ws.Range("L107").Value = 22494
$ws.Range("M107").Value = -2329.5
$ws.Range("N107").Value = -26334
$ws.Range("H132").Value = 3548037.2
$ws.Range("I132").Value = 3548037.2
$ws.Range("K132").Value = 10644111.6
$ws.Range("M132").Value = -10641581.6
$ws.Range("H136").Value = 2587480.5
$ws.Range("I136").Value = 2175643.5
$ws.Range("J136").Value = 2979706.2
$ws.Range("K136").Value = 6526930.5
$ws.Range("L136").Value = 8939118.600000001
$ws.Range("M136").Value = -6524380.5
$ws.Range("N136").Value = -8944218.600000001
